$wb = $excel.ActiveWorkbook

$oldText = "January 30 2026 16.19.47 EST"
$newText = "February 02 2026 12.49.33 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($r + $rowOffset, $c + $colOffset)
            $val = $cell.Value()
            if ($val -ne $null -and $val -is [string] -and $val.Contains($oldText)) {
                $cell.Value = $val.Replace($oldText, $newText)
            }
        }
    }
}
